$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note text introduced by this edit (delta / sort / warning fixes).
# Written in the same order the author typed them so the shared-string
# table grows in the same sequence as the real commit.
$sAddDelta  = "edit to Addelta in utility.R  introduced `nsetdiff(vec_tst,covar) to drop covars`n in onezero vector"
$sStillNeed = "still need check changes correct"
$sSort      = "added to Runmimix.R impdatamergeord[order(impdatamergeord[,"".imp""],`nimpdatamergeord[,idvar]),]`nprior to data output"
$sWarning   = "edit !is.null(method) & (method != ""MAR"")  ) `nalso Warning changed to stop"

# Row 15 -> Addelta/setdiff note (Action) + "still need check changes correct" (Current status)
$ws.Range("I15").Value = $sAddDelta
$ws.Range("I15").WrapText = $true
$ws.Range("I15").VerticalAlignment = -4160

$ws.Range("J15").Value = $sStillNeed
$ws.Range("J15").WrapText = $false
$ws.Range("J15").VerticalAlignment = -4160

# Row 12 -> sort-order fix in Runmimix.R (Action); row grows taller to fit the note
$ws.Rows(12).RowHeight = 75
$ws.Range("I12").Value = $sSort
$ws.Range("I12").WrapText = $true
$ws.Range("I12").VerticalAlignment = -4160

# Row 11 -> warning-changed-to-stop edit (Action)
$ws.Range("I11").Value = $sWarning
$ws.Range("I11").WrapText = $true
$ws.Range("I11").VerticalAlignment = -4160

# Leave the cursor where the author left it.
$ws.Range("I11").Select() | Out-Null
